# Update live market-data columns (H-N) across multiple sheets.
# Values refreshed by the scheduled data-import runner; figures below are the new snapshot.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 1510.0714
$ws.Range("I4").Value = 471.22223
$ws.Range("J4").Value = 3380
$ws.Range("K4").Value = 471.22223
$ws.Range("L4").Value = 3380
$ws.Range("M4").Value = -357.22223
$ws.Range("N4").Value = -3608
$ws.Range("H29").Value = 2754.6365
$ws.Range("H38").Value = 167
$ws.Range("I38").Value = 167
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 501
$ws.Range("L38").Value = 0
$ws.Range("M38").Value = -129
$ws.Range("N38").ClearContents()
$ws.Range("H43").Value = 937.1
$ws.Range("I43").Value = 481
$ws.Range("J43").Value = 1132.5714
$ws.Range("K43").Value = 481
$ws.Range("L43").Value = 1132.5714
$ws.Range("M43").Value = -412
$ws.Range("N43").Value = -1270.5714
$ws.Range("H58").Value = 2250
$ws.Range("I58").Value = 2850
$ws.Range("K58").Value = 8550
$ws.Range("M58").Value = -8400
$ws.Range("H94").Value = 2268.3333
$ws.Range("I94").Value = 2268.3333
$ws.Range("K94").Value = 2268.3333
$ws.Range("M94").Value = -1817.3333
$ws.Range("H98").Value = 2453.5334
$ws.Range("I98").Value = 2121.3
$ws.Range("J98").Value = 3118
$ws.Range("K98").Value = 2121.3
$ws.Range("L98").Value = 3118
$ws.Range("M98").Value = -623.3000000000002
$ws.Range("N98").Value = -6114
$ws.Range("H122").Value = 2453.5334
$ws.Range("I122").Value = 2121.3
$ws.Range("J122").Value = 3118
$ws.Range("K122").Value = 6363.900000000001
$ws.Range("L122").Value = 9354
$ws.Range("M122").Value = -3913.900000000001
$ws.Range("N122").Value = -14254
$ws.Range("H132").Value = 1985.75
$ws.Range("I132").Value = 2115.9092
$ws.Range("K132").Value = 6347.7276
$ws.Range("M132").Value = -3817.7276
$ws.Range("H137").Value = 778213.5600000001
$ws.Range("I137").Value = 3870.1333
$ws.Range("J137").Value = 1193040.4
$ws.Range("K137").Value = 11610.3999
$ws.Range("L137").Value = 3579121.2
$ws.Range("M137").Value = -9060.3999
$ws.Range("N137").Value = -3584221.2
$ws.Range("H138").Value = 4633.5537
$ws.Range("I138").Value = 2192.077
$ws.Range("J138").Value = 5371.6743
$ws.Range("K138").Value = 6576.231000000001
$ws.Range("L138").Value = 16115.0229
$ws.Range("M138").Value = -1436.231000000001
$ws.Range("N138").Value = -26395.0229
$ws.Range("H141").Value = 2959.923
$ws.Range("I141").Value = 1857.88
$ws.Range("J141").Value = 4927.857
$ws.Range("K141").Value = 5573.64
$ws.Range("L141").Value = 14783.571
$ws.Range("M141").Value = -393.6400000000003
$ws.Range("N141").Value = -25143.571

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 20212.793
$ws.Range("I32").Value = 20769.346
$ws.Range("K32").Value = 20769.346
$ws.Range("M32").Value = -20482.346
$ws.Range("H45").Value = 2090.348
$ws.Range("I45").Value = 2049
$ws.Range("K45").Value = 2049
$ws.Range("M45").Value = -1672
$ws.Range("H61").Value = 5446.943
$ws.Range("I61").Value = 4114.3447
$ws.Range("J61").Value = 11887.833
$ws.Range("K61").Value = 4114.3447
$ws.Range("L61").Value = 11887.833
$ws.Range("M61").Value = -3902.3447
$ws.Range("N61").Value = -12311.833
$ws.Range("H63").Value = 3094
$ws.Range("I63").Value = 2805
$ws.Range("J63").Value = 3768.3333
$ws.Range("K63").Value = 2805
$ws.Range("L63").Value = 3768.3333
$ws.Range("M63").Value = -2119
$ws.Range("N63").Value = -5140.3333
$ws.Range("H66").Value = 3094
$ws.Range("I66").Value = 2805
$ws.Range("J66").Value = 3768.3333
$ws.Range("K66").Value = 14025
$ws.Range("L66").Value = 18841.6665
$ws.Range("M66").Value = -10593
$ws.Range("N66").Value = -25705.6665
$ws.Range("H74").Value = 4947.6
$ws.Range("I74").Value = 1687.5
$ws.Range("J74").Value = 17988
$ws.Range("K74").Value = 1687.5
$ws.Range("L74").Value = 17988
$ws.Range("M74").Value = -813.5
$ws.Range("N74").Value = -19736
$ws.Range("H77").Value = 4947.6
$ws.Range("I77").Value = 1687.5
$ws.Range("J77").Value = 17988
$ws.Range("K77").Value = 8437.5
$ws.Range("L77").Value = 89940
$ws.Range("M77").Value = -4069.5
$ws.Range("N77").Value = -98676
$ws.Range("H132").Value = 2474.7346
$ws.Range("I132").Value = 2088.8462
$ws.Range("J132").Value = 3979.7
$ws.Range("K132").Value = 6266.5386
$ws.Range("L132").Value = 11939.1
$ws.Range("M132").Value = -3736.5386
$ws.Range("N132").Value = -16999.1
$ws.Range("H136").Value = 5446.943
$ws.Range("I136").Value = 4114.3447
$ws.Range("J136").Value = 11887.833
$ws.Range("K136").Value = 12343.0341
$ws.Range("L136").Value = 35663.499
$ws.Range("M136").Value = -9793.034099999999
$ws.Range("N136").Value = -40763.499

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H138").Value = 60468.637
$ws.Range("J138").Value = 60468.637
$ws.Range("L138").Value = 60468.637
$ws.Range("N138").Value = -70748.637
$ws.Range("H140").Value = 48959
$ws.Range("J140").Value = 48959
$ws.Range("L140").Value = 48959
$ws.Range("N140").Value = -59319

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 511719.03
$ws.Range("I31").Value = 7112.5483
$ws.Range("J31").Value = 884166.7
$ws.Range("K31").Value = 7112.5483
$ws.Range("L31").Value = 884166.7
$ws.Range("M31").Value = -6817.5483
$ws.Range("N31").Value = -884756.7
$ws.Range("H34").Value = 511719.03
$ws.Range("I34").Value = 7112.5483
$ws.Range("J34").Value = 884166.7
$ws.Range("K34").Value = 7112.5483
$ws.Range("L34").Value = 884166.7
$ws.Range("M34").Value = -6910.5483
$ws.Range("N34").Value = -884570.7
$ws.Range("H80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()
$ws.Range("H107").Value = 1026.75
$ws.Range("I107").Value = 1129.5834
$ws.Range("J107").Value = 923.9167
$ws.Range("K107").Value = 1129.5834
$ws.Range("L107").Value = 923.9167
$ws.Range("M107").Value = 790.4166
$ws.Range("N107").Value = -4763.9167

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 1597.5309
$ws.Range("I68").Value = 1678.4445
$ws.Range("J68").Value = 1532.8
$ws.Range("K68").Value = 5035.333500000001
$ws.Range("L68").Value = 4598.4
$ws.Range("M68").Value = -4224.333500000001
$ws.Range("N68").Value = -6220.4
$ws.Range("H71").Value = 1597.5309
$ws.Range("I71").Value = 1678.4445
$ws.Range("J71").Value = 1532.8
$ws.Range("K71").Value = 15106.0005
$ws.Range("L71").Value = 13795.2
$ws.Range("M71").Value = -11050.0005
$ws.Range("N71").Value = -21907.2
$ws.Range("H92").Value = 304.8
$ws.Range("I92").Value = 312
$ws.Range("J92").Value = 300
$ws.Range("K92").Value = 936
$ws.Range("L92").Value = 900
$ws.Range("M92").Value = 312
$ws.Range("N92").Value = -3396
$ws.Range("H131").Value = 1631.5769
$ws.Range("J131").Value = 1359.5
$ws.Range("L131").Value = 4078.5
$ws.Range("N131").Value = -14158.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3356.5417
$ws.Range("I132").Value = 2394.2222
$ws.Range("J132").Value = 3933.9333
$ws.Range("K132").Value = 7182.6666
$ws.Range("L132").Value = 11801.7999
$ws.Range("M132").Value = -4652.6666
$ws.Range("N132").Value = -16861.7999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 935.6
$ws.Range("J16").Value = 838.5
$ws.Range("L16").Value = 838.5
$ws.Range("N16").Value = -1178.5
$ws.Range("H61").Value = 28491.525
$ws.Range("I61").Value = 27705.193
$ws.Range("J61").Value = 31200
$ws.Range("K61").Value = 27705.193
$ws.Range("L61").Value = 31200
$ws.Range("M61").Value = -27503.193
$ws.Range("N61").Value = -31604
$ws.Range("H113").Value = 28491.525
$ws.Range("I113").Value = 27705.193
$ws.Range("J113").Value = 31200
$ws.Range("K113").Value = 27705.193
$ws.Range("L113").Value = 31200
$ws.Range("M113").Value = -25535.193
$ws.Range("N113").Value = -35540
$ws.Range("H132").Value = 7830.5186
$ws.Range("I132").Value = 9450.799999999999
$ws.Range("J132").Value = 3201.1428
$ws.Range("K132").Value = 28352.4
$ws.Range("L132").Value = 9603.428400000001
$ws.Range("M132").Value = -25822.4
$ws.Range("N132").Value = -14663.4284

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H41").Value = 36500
$ws.Range("J41").Value = 36500
$ws.Range("L41").Value = 36500
$ws.Range("N41").Value = -37280
$ws.Range("H45").Value = 13843
$ws.Range("I45").Value = 13989
$ws.Range("J45").Value = 13806.5
$ws.Range("K45").Value = 13989
$ws.Range("L45").Value = 13806.5
$ws.Range("M45").Value = -13498
$ws.Range("N45").Value = -14788.5
$ws.Range("H107").Value = 4133.4546
$ws.Range("I107").Value = 1154.6
$ws.Range("J107").Value = 6615.8335
$ws.Range("K107").Value = 3463.8
$ws.Range("L107").Value = 19847.5005
$ws.Range("M107").Value = -1543.8
$ws.Range("N107").Value = -23687.5005
